$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 22

# Column A holds a date written as plain text (e.g. "08/02/2025") in every
# other row of this sheet, so write it with a leading apostrophe to stop
# Excel's autocorrect from reinterpreting it as a real date, then strip the
# resulting "quote prefix" formatting so the cell keeps the sheet's default
# (unstyled) look.
$ws.Cells.Item($row, 1).Value = "'06/08/2025"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "Grau"
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = "Deportivo Garcilaso"
$ws.Cells.Item($row, 6).Value = "D"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 1
$ws.Cells.Item($row, 11).Value = 1.65
$ws.Cells.Item($row, 12).Value = 0.89
$ws.Cells.Item($row, 13).Value = 16
$ws.Cells.Item($row, 14).Value = 16
$ws.Cells.Item($row, 15).Value = 5
$ws.Cells.Item($row, 16).Value = 5
